# "Script to populate tables"
# The phone_number column (D) is no longer needed on the nurse sheet, so
# drop it entirely - Excel shifts h_id (E) and the nurse name (F) left
# into D/E for us, which is exactly what the workbook now looks like.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Delete()

# Leave the sheet parked on B1 (top-left data cell) instead of the old
# E16 selection.
[void]$ws.Range("B1").Select()
